# Add 2022-Q3 quarterly data: new sheet + updated summary ("总计") sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Insert a new row for 2022-Q3 at the top of the summary table ---
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()
# carry the bold/bordered style of the index column onto the new row
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$summary.Range("B2:D2").ClearFormats()

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = '2022-Q3'
$summary.Range("C2").Value = 26
$summary.Range("D2").Value = 12.66
$summary.Range("A3").Value = 1
$summary.Range("B3").Value = '2022-Q2'
$summary.Range("C3").Value = 28
$summary.Range("D3").Value = 15.88
$summary.Range("A4").Value = 2
$summary.Range("B4").Value = '2022-Q1'
$summary.Range("C4").Value = 37
$summary.Range("D4").Value = 20.86
$summary.Range("A5").Value = 3
$summary.Range("B5").Value = '2021-Q4'
$summary.Range("C5").Value = 41
$summary.Range("D5").Value = 22.57
$summary.Range("A6").Value = 4
$summary.Range("B6").Value = '2021-Q3'
$summary.Range("C6").Value = 47
$summary.Range("D6").Value = 26.58
$summary.Range("A7").Value = 5
$summary.Range("B7").Value = '2021-Q2'
$summary.Range("C7").Value = 58
$summary.Range("D7").Value = 30
$summary.Range("A8").Value = 6
$summary.Range("B8").Value = '2021-Q1'
$summary.Range("C8").Value = 30
$summary.Range("D8").Value = 18.95
$summary.Range("A9").Value = 7
$summary.Range("B9").Value = '2020-Q4'
$summary.Range("C9").Value = 23
$summary.Range("D9").Value = 17.1

# --- 2. Create the "2022-Q3" sheet (cloned from "2022-Q2" for identical layout/styles) ---
$template = $wb.Worksheets.Item("2022-Q2")
$template.Copy($template)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"
# the template sheet has 28 data rows; 2022-Q3 only has 26, drop the extra two
$newSheet.Range("A28:H29").EntireRow.Delete()

$newSheet.Range("B2").Value = '''001216'
$newSheet.Range("C2").Value = '易方达新收益灵活配置混合 - A'
$newSheet.Range("D2").Value = '''39.46'
$newSheet.Range("E2").Value = '''77.61'
$newSheet.Range("F2").Value = '''5.29'
$newSheet.Range("G2").Value = '''2.0874'
$newSheet.Range("H2").Value = 4
$newSheet.Range("B3").Value = '''005395'
$newSheet.Range("C3").Value = '泓德臻远回报灵活配置混合'
$newSheet.Range("D3").Value = '''29.76'
$newSheet.Range("E3").Value = '''93.35'
$newSheet.Range("F3").Value = '''6.95'
$newSheet.Range("G3").Value = '''2.0683'
$newSheet.Range("H3").Value = 3
$newSheet.Range("B4").Value = '''001500'
$newSheet.Range("C4").Value = '泓德远见回报混合'
$newSheet.Range("D4").Value = '''21.21'
$newSheet.Range("E4").Value = '''92.75'
$newSheet.Range("F4").Value = '''7.24'
$newSheet.Range("G4").Value = '''1.5356'
$newSheet.Range("H4").Value = 6
$newSheet.Range("B5").Value = '''519692'
$newSheet.Range("C5").Value = '交银成长混合A'
$newSheet.Range("D5").Value = '''23.45'
$newSheet.Range("E5").Value = '''76.71'
$newSheet.Range("F5").Value = '''6.02'
$newSheet.Range("G5").Value = '''1.4117'
$newSheet.Range("H5").Value = 7
$newSheet.Range("B6").Value = '''001603'
$newSheet.Range("C6").Value = '易方达安盈回报混合'
$newSheet.Range("D6").Value = '''25.59'
$newSheet.Range("E6").Value = '''32.33'
$newSheet.Range("F6").Value = '''3.93'
$newSheet.Range("G6").Value = '''1.0057'
$newSheet.Range("H6").Value = 2
$newSheet.Range("B7").Value = '''519694'
$newSheet.Range("C7").Value = '交银蓝筹混合'
$newSheet.Range("D7").Value = '''16.10'
$newSheet.Range("E7").Value = '''78.09'
$newSheet.Range("F7").Value = '''6.20'
$newSheet.Range("G7").Value = '''0.9982'
$newSheet.Range("H7").Value = 7
$newSheet.Range("B8").Value = '''004965'
$newSheet.Range("C8").Value = '泓德致远混合A'
$newSheet.Range("D8").Value = '''16.19'
$newSheet.Range("E8").Value = '''46.90'
$newSheet.Range("F8").Value = '''4.69'
$newSheet.Range("G8").Value = '''0.7593'
$newSheet.Range("H8").Value = 4
$newSheet.Range("B9").Value = '''001217'
$newSheet.Range("C9").Value = '易方达新收益灵活配置混合 - C'
$newSheet.Range("D9").Value = '''14.06'
$newSheet.Range("E9").Value = '''77.61'
$newSheet.Range("F9").Value = '''5.29'
$newSheet.Range("G9").Value = '''0.7438'
$newSheet.Range("H9").Value = 4
$newSheet.Range("B10").Value = '''009812'
$newSheet.Range("C10").Value = '易方达悦兴一年持有期混合A'
$newSheet.Range("D10").Value = '''61.64'
$newSheet.Range("E10").Value = '''20.74'
$newSheet.Range("F10").Value = '''0.77'
$newSheet.Range("G10").Value = '''0.4746'
$newSheet.Range("H10").Value = 10
$newSheet.Range("B11").Value = '''004848'
$newSheet.Range("C11").Value = '中欧睿泓定期开放灵活配置混合'
$newSheet.Range("D11").Value = '''17.10'
$newSheet.Range("E11").Value = '''59.07'
$newSheet.Range("F11").Value = '''2.27'
$newSheet.Range("G11").Value = '''0.3882'
$newSheet.Range("H11").Value = 8
$newSheet.Range("B12").Value = '''011837'
$newSheet.Range("C12").Value = '鹏扬中国优质成长混合A'
$newSheet.Range("D12").Value = '''7.87'
$newSheet.Range("E12").Value = '''88.18'
$newSheet.Range("F12").Value = '''3.09'
$newSheet.Range("G12").Value = '''0.2432'
$newSheet.Range("H12").Value = 9
$newSheet.Range("B13").Value = '''001319'
$newSheet.Range("C13").Value = '农银汇理信息传媒主题股票'
$newSheet.Range("D13").Value = '''2.62'
$newSheet.Range("E13").Value = '''80.55'
$newSheet.Range("F13").Value = '''7.91'
$newSheet.Range("G13").Value = '''0.2072'
$newSheet.Range("H13").Value = 5
$newSheet.Range("B14").Value = '''009813'
$newSheet.Range("C14").Value = '易方达悦兴一年持有期混合C'
$newSheet.Range("D14").Value = '''17.96'
$newSheet.Range("E14").Value = '''20.74'
$newSheet.Range("F14").Value = '''0.77'
$newSheet.Range("G14").Value = '''0.1383'
$newSheet.Range("H14").Value = 10
$newSheet.Range("B15").Value = '''004966'
$newSheet.Range("C15").Value = '泓德致远混合C'
$newSheet.Range("D15").Value = '''2.54'
$newSheet.Range("E15").Value = '''46.90'
$newSheet.Range("F15").Value = '''4.69'
$newSheet.Range("G15").Value = '''0.1191'
$newSheet.Range("H15").Value = 4
$newSheet.Range("B16").Value = '''005642'
$newSheet.Range("C16").Value = '鹏扬景升灵活配置混合A'
$newSheet.Range("D16").Value = '''3.35'
$newSheet.Range("E16").Value = '''87.37'
$newSheet.Range("F16").Value = '''3.39'
$newSheet.Range("G16").Value = '''0.1136'
$newSheet.Range("H16").Value = 9
$newSheet.Range("B17").Value = '''005664'
$newSheet.Range("C17").Value = '鹏扬景欣混合A'
$newSheet.Range("D17").Value = '''6.93'
$newSheet.Range("E17").Value = '''26.14'
$newSheet.Range("F17").Value = '''1.13'
$newSheet.Range("G17").Value = '''0.0783'
$newSheet.Range("H17").Value = 9
$newSheet.Range("B18").Value = '''006977'
$newSheet.Range("C18").Value = '农银汇理海棠三年定期开放混合'
$newSheet.Range("D18").Value = '''1.40'
$newSheet.Range("E18").Value = '''80.33'
$newSheet.Range("F18").Value = '''5.17'
$newSheet.Range("G18").Value = '''0.0724'
$newSheet.Range("H18").Value = 9
$newSheet.Range("B19").Value = '''004341'
$newSheet.Range("C19").Value = '农银汇理尖端科技灵活配置混合'
$newSheet.Range("D19").Value = '''1.10'
$newSheet.Range("E19").Value = '''85.97'
$newSheet.Range("F19").Value = '''5.11'
$newSheet.Range("G19").Value = '''0.0562'
$newSheet.Range("H19").Value = 9
$newSheet.Range("B20").Value = '''001060'
$newSheet.Range("C20").Value = '前海开源高端装备制造灵活配置混合'
$newSheet.Range("D20").Value = '''0.90'
$newSheet.Range("E20").Value = '''80.76'
$newSheet.Range("F20").Value = '''4.65'
$newSheet.Range("G20").Value = '''0.0418'
$newSheet.Range("H20").Value = 7
$newSheet.Range("B21").Value = '''011838'
$newSheet.Range("C21").Value = '鹏扬中国优质成长混合C'
$newSheet.Range("D21").Value = '''1.29'
$newSheet.Range("E21").Value = '''88.18'
$newSheet.Range("F21").Value = '''3.09'
$newSheet.Range("G21").Value = '''0.0399'
$newSheet.Range("H21").Value = 9
$newSheet.Range("B22").Value = '''008499'
$newSheet.Range("C22").Value = '鹏扬景科混合A'
$newSheet.Range("D22").Value = '''1.69'
$newSheet.Range("E22").Value = '''33.87'
$newSheet.Range("F22").Value = '''1.38'
$newSheet.Range("G22").Value = '''0.0233'
$newSheet.Range("H22").Value = 7
$newSheet.Range("B23").Value = '''005643'
$newSheet.Range("C23").Value = '鹏扬景升灵活配置混合C'
$newSheet.Range("D23").Value = '''0.67'
$newSheet.Range("E23").Value = '''87.37'
$newSheet.Range("F23").Value = '''3.39'
$newSheet.Range("G23").Value = '''0.0227'
$newSheet.Range("H23").Value = 9
$newSheet.Range("B24").Value = '''005665'
$newSheet.Range("C24").Value = '鹏扬景欣混合C'
$newSheet.Range("D24").Value = '''1.40'
$newSheet.Range("E24").Value = '''26.14'
$newSheet.Range("F24").Value = '''1.13'
$newSheet.Range("G24").Value = '''0.0158'
$newSheet.Range("H24").Value = 9
$newSheet.Range("B25").Value = '''008500'
$newSheet.Range("C25").Value = '鹏扬景科混合C'
$newSheet.Range("D25").Value = '''0.72'
$newSheet.Range("E25").Value = '''33.87'
$newSheet.Range("F25").Value = '''1.38'
$newSheet.Range("G25").Value = '''0.0099'
$newSheet.Range("H25").Value = 7
$newSheet.Range("B26").Value = '''960016'
$newSheet.Range("C26").Value = '交银成长混合H'
$newSheet.Range("D26").Value = '''0.16'
$newSheet.Range("E26").Value = '''76.71'
$newSheet.Range("F26").Value = '''6.02'
$newSheet.Range("G26").Value = '''0.0096'
$newSheet.Range("H26").Value = 7
$newSheet.Range("B27").Value = '''001708'
$newSheet.Range("C27").Value = '东兴改革精选灵活配置混合'
$newSheet.Range("D27").Value = '''0.03'
$newSheet.Range("E27").Value = '''90.16'
$newSheet.Range("F27").Value = '''2.74'
$newSheet.Range("G27").Value = '''0.0008'
$newSheet.Range("H27").Value = 9
# drop the quote-prefix style injected by the leading apostrophes above
$newSheet.Range("B2:G27").Style = "Normal"

$summary.Select()
